$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Re-assert the table's preferred width (100%) so it is re-serialized as an
# integer percentage (5000) instead of the legacy floating-point form.
$tbl.PreferredWidthType = 2
$tbl.PreferredWidth = 250

# Mark the first row as a repeating header row (w:trPr/w:tblHeader)
$tbl.Rows.Item(1).HeadingFormat = $true

# Locate the row containing the Oreo filling confidence-interval answer and
# update the numeric bounds (2.808/2.988 -> 2.535/3.165).
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cell = $tbl.Cell($r, 3)
    if ($cell.Range.Text -like "*2.808*2.988*Oreo filling*") {
        $para = $cell.Range.Paragraphs.Item(1)
        $start = $para.Range.Start

        $run1 = $d.Range($start, $start + 78)
        $run1.Text = "(2.535 , 3.165) We are 90% confident that the true mean weight of Oreo filling"

        $run3 = $d.Range($start + 79, $start + 111)
        $run3.Text = "is between 2.535 and 3.165 grams"

        break
    }
}
